$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Format the Price column as Text first so numeric-looking strings
# (e.g. "1.007") are stored verbatim instead of being parsed as numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "27.015.48"
$ws.Range("E2").Value = "  -1.08%  "
$ws.Range("D3").Value = "1.826.67"
$ws.Range("E3").Value = "  -0.38%  "
$ws.Range("D4").Value = "1.007"
$ws.Range("E4").Value = "  -0.38%  "
$ws.Range("D5").Value = "310.39"
$ws.Range("E5").Value = "  -1.30%  "
$ws.Range("D6").Value = "1.005"
$ws.Range("E6").Value = "  -0.41%  "
$ws.Range("D7").Value = "0.4629"
$ws.Range("E7").Value = "  -2.26%  "
$ws.Range("D8").Value = "0.3731"
$ws.Range("E8").Value = "  +1.20%  "
$ws.Range("D9").Value = "0.07249"
$ws.Range("E9").Value = "  -2.67%  "
$ws.Range("D10").Value = "0.8629"
$ws.Range("E10").Value = "  -2.68%  "
$ws.Range("D11").Value = "19.94"
$ws.Range("E11").Value = "  -2.72%  "
$ws.Range("D12").Value = "0.07795"
$ws.Range("E12").Value = "  +6.33%  "
$ws.Range("E13").Value = "  -2.34%  "
$ws.Range("D14").Value = "5.342"
$ws.Range("E14").Value = "  -1.64%  "
$ws.Range("D15").Value = "6.533"
$ws.Range("E15").Value = "  -0.52%  "
$ws.Range("D16").Value = "91.80"
$ws.Range("E16").Value = "  -2.46%  "
$ws.Range("D17").Value = "1.007"
$ws.Range("E17").Value = "  -0.11%  "
$ws.Range("D18").Value = "0.000008686"
$ws.Range("E18").Value = "  -1.28%  "
$ws.Range("E19").Value = "  -0.36%  "
$ws.Range("D20").Value = "27.141.86"
$ws.Range("E20").Value = "  -1.61%  "
$ws.Range("D21").Value = "14.54"
$ws.Range("E21").Value = "  -1.67%  "
$ws.Range("D22").Value = "5.154"
$ws.Range("E22").Value = "  -2.54%  "
$ws.Range("D23").Value = "10.56"
$ws.Range("E23").Value = "  -1.12%  "
$ws.Range("D24").Value = "2.077.20"
$ws.Range("E24").Value = "  -1.23%  "
$ws.Range("D25").Value = "153.01"
$ws.Range("E25").Value = "  +0.75%  "
$ws.Range("D26").Value = "1.838"
$ws.Range("E26").Value = "  -3.07%  "
$ws.Range("D27").Value = "18.18"
$ws.Range("E27").Value = "  -2.54%  "
$ws.Range("D28").Value = "2.087"
$ws.Range("E28").Value = "  -2.59%  "
$ws.Range("D29").Value = "5.120"
$ws.Range("E29").Value = "  -2.16%  "
$ws.Range("D30").Value = "115.36"
$ws.Range("E30").Value = "  -1.71%  "
$ws.Range("D31").Value = "0.08844"
$ws.Range("E31").Value = "  -1.67%  "
$ws.Range("D32").Value = "2.965"
$ws.Range("E32").Value = "  +0.51%  "
$ws.Range("D33").Value = "0.7280"
$ws.Range("E33").Value = "  -3.01%  "
$ws.Range("D34").Value = "4.435"
$ws.Range("E34").Value = "  -2.48%  "
$ws.Range("D35").Value = "1.134"
$ws.Range("E35").Value = "  -3.48%  "
$ws.Range("D36").Value = "2.491"
$ws.Range("E36").Value = "  +2.34%  "
$ws.Range("E37").Value = "  -1.46%  "
$ws.Range("D38").Value = "0.01944"
$ws.Range("E38").Value = "  -0.64%  "
$ws.Range("D39").Value = "0.05231"
$ws.Range("E39").Value = "  -2.04%  "
$ws.Range("D40").Value = "2.933"
$ws.Range("E40").Value = "  -1.02%  "
$ws.Range("D41").Value = "7.225"
$ws.Range("E41").Value = "  -0.34%  "
$ws.Range("D42").Value = "0.5157"
$ws.Range("E42").Value = "  -2.56%  "
$ws.Range("D43").Value = "0.1627"
$ws.Range("E43").Value = "  -1.93%  "
$ws.Range("D44").Value = "0.8565"
$ws.Range("E44").Value = "  -15.27%  "
$ws.Range("D45").Value = "8.194"
$ws.Range("E45").Value = "  -3.56%  "
$ws.Range("D46").Value = "0.4813"
$ws.Range("E46").Value = "  -2.29%  "
$ws.Range("D47").Value = "1.005"
$ws.Range("E47").Value = "  -0.52%  "
$ws.Range("D48").Value = "10.20"
$ws.Range("E48").Value = "  -2.79%  "
$ws.Range("D49").Value = "102.68"
$ws.Range("E49").Value = "  -2.44%  "
$ws.Range("E50").Value = "  -0.63%  "
$ws.Range("D51").Value = "1.619"
$ws.Range("E51").Value = "  -2.83%  "

# Restore the original (default) cell style now that the text is locked in.
$ws.Range("D2:D51").Style = "Normal"
